# Applies the odds updates described in the commit diff for row 3 and row 4
# of Sheet1 ("Jogos_da_Semana_FlashScore_2024-10-10.xlsx").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 updates ---
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 8
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 29
$ws.Range("AK3").Value = 67
$ws.Range("AM3").Value = 51
$ws.Range("AS3").Value = 201

# --- Row 4 updates ---
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 2.65
$ws.Range("P4").Value = 3.55
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 2.02
$ws.Range("W4").Value = 11.75
$ws.Range("X4").Value = 18
$ws.Range("Z4").Value = 40
$ws.Range("AA4").Value = 24
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 6.9
$ws.Range("AE4").Value = 12
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 11.75
$ws.Range("AM4").Value = 22
$ws.Range("AO4").Value = 16
$ws.Range("AP4").Value = 21
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.95

$wb.Save()
